$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# B2: Cases query - remove the Cohort clause (no longer join cohort)
$b2Text = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE s.clinical_study_designation IN ['UBC01']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@
$ws.Range("B2").Value = $b2Text

# C2/C3/C4: StatQuery - study designation UBC01 -> UBC02
$statText = @'
MATCH (s:study)
WHERE
	s.clinical_study_designation IN ['UBC02']
OPTIONAL MATCH (s)<-[:member_of]-(c:case)
OPTIONAL MATCH (c)<-[:of_case]-(samp:sample)<-[:of_sample]-(f:file)
RETURN 
	count(DISTINCT(f)) as number_of_files , 
	count(DISTINCT(samp)) as number_of_sample , 
	count(DISTINCT(c.case_id)) as number_of_cases , 
	count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@
$ws.Range("C2").Value = $statText
$ws.Range("C3").Value = $statText
$ws.Range("C4").Value = $statText

# Row heights (re-flowed wrap-text heights)
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 230.4
$ws.Rows.Item(4).RowHeight = 230.4

# Column widths (minor re-measurement)
$ws.Columns.Item(1).ColumnWidth = 10.88671875
$ws.Columns.Item(2).ColumnWidth = 75.77734375
$ws.Columns.Item(3).ColumnWidth = 75.77734375
$ws.Columns.Item(4).ColumnWidth = 70.21875
$ws.Columns.Item(5).ColumnWidth = 39.77734375

# View state: selection moved from H4 to B4, scrolled to A4
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select()
